$wb = $excel.ActiveWorkbook

# The "Global" sheet is duplicated to create a new parameter sheet named
# "Sheet, with comma", inserted right after "Global" (before "MissingParam").
$global = $wb.Worksheets.Item("Global")
$global.Copy($null, $global)

# The copy becomes the active sheet and is placed immediately after Global.
$newSheet = $wb.ActiveSheet
$newSheet.Name = "Sheet, with comma"
$newSheet.Range("A30").Select()

# Update the selection remembered on the Global sheet (whole-sheet selection).
$global.Cells.Select()

# Keep the new sheet as the active/selected tab.
$newSheet.Activate()
$newSheet.Range("A30").Select()
